$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row: rename column B header from "Name" to "Asset Name"
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "Asset Name"

# ---------------------------------------------------------------------------
# New sample data row (row 2) - import/export example asset
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "Build Permanent-Housing"
$ws.Range("C2").Value = 3304
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 20
$ws.Range("F2").Value = 1000
$ws.Range("G2").Value = 1000

# ---------------------------------------------------------------------------
# Apply the built-in "Neutral" cell style (yellow fill / brown text) to the
# new data row, matching the exact colors Excel uses for that style.
# ---------------------------------------------------------------------------
$ws.Range("A2:G2").Style = "Neutral"

$neutral = $wb.Styles.Item("Neutral")
$neutral.Font.Color = 26012
$neutral.Interior.Color = 10284031

# Re-establish the vertical-center alignment (inherited from the template's
# original body-row formatting) for the first three columns of the new row.
$ws.Range("A2:C2").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Column widths - widen to fit the newly-entered sample data / headers.
# ---------------------------------------------------------------------------
$ws.Columns("B:B").ColumnWidth = 20.944010416666668
$ws.Columns("E:E").ColumnWidth = 15.053385416666666
$ws.Columns("F:F").ColumnWidth = 13.276041666666666
$ws.Columns("G:G").ColumnWidth = 14.721354166666666

# ---------------------------------------------------------------------------
# Update the active selection shown when the workbook is reopened.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("D10").Select()
